$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 13:12"

# Full refreshed province/city table: city name + Casos totales, Casos activos,
# Recuperados, Muertes for rows 4-64 (data updated & re-sorted by Casos totales desc)
$data = @(
    @(4, 'Madrid', 19243, 5044, 11787, 2412),
    @(5, 'Cataluña', 12940, 2384, 9676, 880),
    @(6, 'Castilla-La Mancha', 2780, 71, 2446, 263),
    @(7, 'Bizkaia/Vizcaya', 2186, 814, 1793, 79),
    @(8, 'Valencia/Valencia', 2027, 50, 1889, 88),
    @(9, 'Navarra', 1641, 70, 1513, 58),
    @(10, 'Araba/Alava', 1610, 814, 1207, 103),
    @(11, 'La Rioja', 1236, 62, 1119, 55),
    @(12, 'Ciudad Real', 1147, 95, 1050, 89),
    @(13, 'Alacant/Alicante', 1093, 19, 983, 91),
    @(14, 'A Coruña', 1049, 67, 1000, 27),
    @(15, 'Zaragoza', 1045, 68, 928, 49),
    @(16, 'Toledo', 965, 95, 860, 78),
    @(17, 'Aragon', 907, 29, 838, 40),
    @(18, 'Malaga', 905, 52, 811, 42),
    @(19, 'Asturias', 900, 52, 819, 29),
    @(20, 'Gran Canaria', 878, 18, 818, 27),
    @(21, 'Cantabria', 810, 19, 770, 21),
    @(22, 'Gipuzkoa/Guipuzcoa', 805, 814, 580, 25),
    @(23, 'Pontevedra', 795, 67, 769, 8),
    @(24, 'Murcia', 687, 12, 660, 15),
    @(25, 'Albacete', 666, 95, 592, 66),
    @(26, 'Granada', 648, 1, 615, 32),
    @(27, 'Caceres', 640, 4, 587, 35),
    @(28, 'Salamanca', 629, 73, 502, 54),
    @(29, 'Sevilla', 627, 8, 595, 24),
    @(30, 'Valladolid', 598, 50, 521, 27),
    @(31, 'Tenerife', 539, 15, 519, 27),
    @(32, 'Burgos', 530, 67, 431, 32),
    @(33, 'Leon', 516, 48, 422, 46),
    @(34, 'Guadalajara', 428, 95, 370, 56),
    @(35, 'Castello/Castellon', 412, 4, 389, 19),
    @(36, 'Segovia', 361, 62, 262, 37),
    @(37, 'Cordoba', 349, 0, 342, 7),
    @(38, 'Jaen', 349, 5, 329, 15),
    @(39, 'Soria', 339, 32, 291, 16),
    @(40, 'Badajoz', 329, 20, 316, 4),
    @(41, 'Ourense', 321, 67, 302, 5),
    @(42, 'Cadiz', 299, 4, 290, 5),
    @(43, 'Avila', 270, 53, 190, 27),
    @(44, 'Mallorca', 210, 18, 194, 12),
    @(45, 'Cuenca', 177, 95, 142, 27),
    @(46, 'Lugo', 157, 67, 140, 4),
    @(47, 'Huesca', 150, 10, 136, 4),
    @(48, 'Teruel', 143, 9, 129, 5),
    @(49, 'Palencia', 139, 14, 120, 5),
    @(50, 'Almeria', 134, 5, 122, 7),
    @(51, 'Zamora', 106, 24, 73, 9),
    @(52, 'Huelva', 95, 2, 91, 2),
    @(53, 'Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena', 58, 0, 58, 3),
    @(54, 'Melilla', 39, 0, 38, 1),
    @(55, 'La Palma', 33, 15, 33, 27),
    @(56, 'Ibiza', 21, 18, 20, 1),
    @(57, 'Fuerteventura', 20, 15, 20, 27),
    @(58, 'Lanzarote', 17, 15, 17, 27),
    @(59, 'Ceuta', 16, 0, 16, 0),
    @(60, 'Menorca', 15, 18, 13, 0),
    @(61, 'Arroyo de la Luz', 7, 0, 7, 0),
    @(62, 'La Gomera', 4, 15, 2, 27),
    @(63, 'El Hierro', 3, 15, 3, 27),
    @(64, 'Formentera', 0, 10, 0, 8)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
